# Update Sheet1 with new travel-insurance plan data (bat files added and notification handled)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "Royal Sundaram"
$ws.Range("B2").Value = "₹874GST included"

$ws.Range("A3").Value = "Travel Shield Single Trip"
$ws.Range("B3").Value = "₹1,220GST included"

$ws.Range("A4").Value = "Care Health"
$ws.Range("B4").Value = "₹1,333GST included"
